$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("DONE", "TECHM.NS",   "HIGH", 1347, "Pending", "2024-06-11 06:44:41", "New data added"),
    @("ADD",  "TECHM.NS",   "HIGH", 1375, "Pending", "2024-06-11 06:44:41", "New data added"),
    @("ADD",  "IRCON.NS",   "HIGH", 267,  "Pending", "2024-06-11 06:44:41", "New data added"),
    @("ADD",  "IRFC.NS",    "HIGH", 183,  "Pending", "2024-06-11 06:44:41", "New data added"),
    @("ADD",  "RAILTEL.NS", "HIGH", 420,  "Pending", "2024-06-11 06:44:41", "New data added"),
    @("ADD",  "RAILTEL.NS", "HIGH", 415,  "Pending", "2024-06-11 06:44:41", "New data added")
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
